$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.354.89"
$ws.Range("E2").Value = "  +0.50%  "

$ws.Range("D3").Value = "2.239.79"

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").Value = "'245.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.40%  "

$ws.Range("D6").Value = "'0.630"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.51%  "

$ws.Range("D7").Value = "'74.29"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.67%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("E9").Value = "  -0.44%  "

$ws.Range("D10").Value = "'43.56"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.54%  "

$ws.Range("D11").Value = "'0.0957"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.33%  "

$ws.Range("D12").Value = "'7.13"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.56%  "

$ws.Range("E13").Value = "  +0.81%  "

$ws.Range("D14").Value = "'14.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.28%  "

$ws.Range("D15").Value = "'0.854"
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").Value = "2.246.40"
$ws.Range("E16").Value = "  +0.38%  "

$ws.Range("D17").Value = "42.276.76"
$ws.Range("E17").Value = "  +0.58%  "

$ws.Range("E18").Value = "  +11.65%  "

$ws.Range("E19").Value = "  +1.14%  "

$ws.Range("D20").Value = "'72.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.32%  "

$ws.Range("D21").Value = "'10.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +41.04%  "

$ws.Range("D22").Value = "'231.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("E23").Value = "  -4.69%  "

$ws.Range("D24").Value = "'11.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.41%  "

$ws.Range("E25").Value = "  +0.10%  "

$ws.Range("E26").Value = "  +0.49%  "

$ws.Range("E27").Value = "  +0.41%  "

$ws.Range("D28").Value = "'2.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.12%  "

$ws.Range("D29").Value = "'166.99"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.43%  "

$ws.Range("E30").Value = "  +1.59%  "

$ws.Range("D31").Value = "'5.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +20.23%  "

$ws.Range("D32").Value = "'0.0813"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.75%  "

$ws.Range("D33").Value = "'0.118"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.62%  "

$ws.Range("D34").Value = "'29.96"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -9.63%  "

$ws.Range("E35").Value = "  -0.10%  "

$ws.Range("D36").Value = "'4.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.43%  "

$ws.Range("E37").Value = "  +2.57%  "

$ws.Range("D38").Value = "'13.36"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.02%  "

$ws.Range("E39").Value = "  -0.29%  "

$ws.Range("E40").Value = "  -3.95%  "

$ws.Range("D41").Value = "'63.53"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.79%  "

$ws.Range("D42").Value = "'0.202"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.57%  "

$ws.Range("D43").Value = "'8.85"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.16%  "

$ws.Range("D44").Value = "'105.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.99%  "

$ws.Range("D45").Value = "'0.102"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.52%  "

$ws.Range("E46").Value = "  -0.18%  "

$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "'1.14"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.53%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'2.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.64%  "

$ws.Range("E49").Value = "  +1.14%  "

$ws.Range("D50").Value = "'2.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.47%  "

$ws.Range("D51").Value = "'4.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.81%  "
